$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "C:\Users\Veeraraju_elluru\Downloads\image_tagging_app\image_tagging_app\uploads\a1.jpg"
$ws.Range("B2").Value = "flower, blue"

$ws.Range("A3").Value = "C:\Users\Veeraraju_elluru\Downloads\image_tagging_app\image_tagging_app\uploads\a3.jpg"
$ws.Range("B3").Value = "blue, pot, flower, plants"
